{"js": "// Add \"Bootstrap\" to the platform/framework requirements bullet, and move\n// the `_GoBack` bookmark (an artifact of the author's last cursor position)\n// from the old \"Registration and Authentication\" bullet to the new edit\n// point in the platform bullet.\n\n// 1) Remove the stray `_GoBack` bookmark currently sitting inside\n//    \"Registration and Authentic|ation\" - Word will drop a fresh one at the\n//    new last-edit location instead.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Re-normalize that paragraph's text into a single run (Word merges runs\n//    with identical formatting back together once the bookmark split is\n//    gone).\nconst regResults = context.document.body.search(\"Registration and Authentication\", { matchCase: true });\nregResults.load(\"items/text\");\nawait context.sync();\nif (regResults.items.length > 0) {\n  regResults.items[0].insertText(\"Registration and Authentication\", \"Replace\");\n  await context.sync();\n}\n\n// 3) Insert \" and Bootstrap\" right after \"(MVC)\" in the platform bullet.\nconst mvcResults = context.document.body.search(\"Site to be developed using PHP Model View Control (MVC)\", { matchCase: true });\nmvcResults.load(\"items/text\");\nawait context.sync();\nif (mvcResults.items.length > 0) {\n  mvcResults.items[0].insertText(\" and Bootstrap\", \"After\");\n  await context.sync();\n}\n\n// 4) Re-insert the `_GoBack` bookmark immediately before \" with My SQL\n//    Backend\" - i.e. at the new last-edit point.\nconst tailResults = context.document.body.search(\" with My SQL Backend\", { matchCase: true });\ntailResults.load(\"items/text\");\nawait context.sync();\nif (tailResults.items.length > 0) {\n  const insertionPoint = tailResults.items[0].getRange(\"Start\");\n  insertionPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Add \"Bootstrap\" to the platform/framework requirements bullet, and move\n# the `_GoBack` bookmark (an artifact of the author's last cursor position)\n# from the old \"Registration and Authentication\" bullet to the new edit\n# point in the platform bullet.\n\n$d = $word.ActiveDocument\n\n# Word constants used below.\n$wdReplaceOne   = 2\n$wdCollapseStart = 1\n$wdCollapseEnd   = 0\n\n# 1) Remove the stray `_GoBack` bookmark currently sitting inside\n#    \"Registration and Authentic|ation\".\n$d.Bookmarks(\"_GoBack\").Delete()\n\n# 2) Re-normalize that paragraph's text into a single clean run (use\n#    Find/Replace so the surrounding character formatting is preserved).\n$regRange = $d.Content\n$regFound = $regRange.Find.Execute(\n    \"Registration and Authentication\",\n    $true, $false, $false, $false, $false, $true, 1, $false,\n    \"Registration and Authentication\", $wdReplaceOne)\nWrite-Output \"registration replace found=$regFound\"\n\n# 3) Insert \" and Bootstrap\" right after \"(MVC)\" in the platform bullet.\n$mvcRange = $d.Content\n$mvcFound = $mvcRange.Find.Execute(\"Site to be developed using PHP Model View Control (MVC)\", $true)\nWrite-Output \"mvc found=$mvcFound\"\nif ($mvcFound) {\n    $mvcRange.Collapse($wdCollapseEnd)\n    $mvcRange.InsertAfter(\" and Bootstrap\")\n}\n\n# 4) Re-insert the `_GoBack` bookmark immediately before \" with My SQL\n#    Backend\" - i.e. at the new last-edit point.\n$tailRange = $d.Content\n$tailFound = $tailRange.Find.Execute(\" with My SQL Backend\", $true)\nWrite-Output \"tail found=$tailFound\"\nif ($tailFound) {\n    $tailRange.Collapse($wdCollapseStart)\n    $d.Bookmarks.Add(\"_GoBack\", $tailRange)\n}\n"}
